# Update "想去人数" (attendance) figures on the 展览 (Exhibition) and
# 全部类型 (All Types) sheets to the freshly scraped counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 706
$ws1.Range("F3").Value = 51
$ws1.Range("F4").Value = 544
$ws1.Range("F9").Value = 4507
$ws1.Range("F10").Value = 4379
$ws1.Range("F11").Value = 9

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 706
$ws4.Range("F3").Value = 51
$ws4.Range("F4").Value = 544
$ws4.Range("F9").Value = 4508
$ws4.Range("F10").Value = 4379
$ws4.Range("F11").Value = 9
